# Add a "home timezone" row to the Tournament info table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")

# Make Tournament the active sheet/tab (matches workbookView losing its
# explicit activeTab, i.e. reverting to the first sheet).
$ws.Activate()

# Insert a new row above the current row 5 ("venue.1" / Atlanta), shifting
# all the venue rows down by one.
$ws.Rows("5:5").Insert()

# Populate the new row with the timezone key/value; other columns (es, it,
# fr, de, nl, ja, fa) are left blank for this row.
$ws.Range("A5").Value = "timezone"
$ws.Range("B5").Value = "America/New_York"

# Grow the "tournament" table to include the newly inserted row.
$lo = $ws.ListObjects.Item("tournament")
$lo.Resize($ws.Range("A1:I16"))

# Reflect the row-5 selection left by the edit.
$ws.Range("A5:XFD5").Select()
